$d = $word.ActiveDocument

# Anchor on the "Is there a correlation..." question, then walk forward:
#   heading -> scatter-plot picture paragraph -> first empty paragraph (to edit)
#           -> second empty paragraph (to delete) -> third empty paragraph (untouched)
$heading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Is there a correlation between number of all-inclusive hotels and score?*") {
        $heading = $p
        break
    }
}

if ($heading -eq $null) {
    throw "Could not find the correlation question paragraph"
}

$picturePara = $heading.Next()
$target = $picturePara.Next()
$emptyFollower = $target.Next()

if ($target.Range.Text -ne "`r" -or $emptyFollower.Range.Text -ne "`r") {
    throw "Unexpected document shape around the correlation answer; aborting"
}

# Build the replacement run as literal OOXML so the inserted run carries
# exactly the same run formatting (w:lang="en-US") used throughout the
# document, then drop it right before the (still empty) paragraph mark.
$insertPoint = $d.Range($target.Range.Start, $target.Range.Start)
$text = "With scatter plot we can see if it exists a correlation between two variables, in our case feedback and the number of all inclusive hotels. We can see it does not exist a correlation. If it existed, we would see a line. "
$ooxml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">' + $text + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($ooxml)

# Justify the paragraph that now holds the new sentence.
$target.Format.Alignment = "wdAlignParagraphJustify"

# Remove the second (now redundant) empty paragraph.
$emptyFollower.Range.Delete()

Write-Output "ok"
